$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.722.75"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "'1.533.18"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'205.63"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'21.38"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'1.749.24"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'1.529.43"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "'26.707.20"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'61.06"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'212.02"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "'7.20"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'9.07"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "'151.88"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "'14.80"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").Value = "'1.10"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "'0.0453"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "'1.360.86"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'0.944"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.798"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.71"
$ws.Range("E41").Value = "  +6.20%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'0.993"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'62.51"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'1.665.01"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'85.33"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").Value = "'0.0₇0969"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "'0.0943"
$ws.Range("E51").Value = "  -0.49%  "
